$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 1.38
$ws.Range("J2").Value = 5.4
$ws.Range("K2").Value = 5.5
$ws.Range("T2").Value = 2.62
$ws.Range("W2").Value = 3.65
$ws.Range("Y2").Value = 29
$ws.Range("AB3").Value = 12
$ws.Range("H3").Value = 2.74
$ws.Range("J3").Value = 3.35
$ws.Range("K3").Value = 3.4
$ws.Range("L3").Value = 1.43
$ws.Range("M3").Value = 1.08
$ws.Range("T3").Value = 1.78
$ws.Range("U3").Value = 2.22
$ws.Range("V3").Value = 1.56
$ws.Range("AD4").Value = 29
$ws.Range("AG4").Value = 10
$ws.Range("Q4").Value = 1.46
$ws.Range("R4").Value = 1.85
$ws.Range("S4").Value = 2.14
$ws.Range("T4").Value = 1.68
$ws.Range("U4").Value = 2.42
$ws.Range("AB5").Value = 12
$ws.Range("AC5").Value = 8.4
$ws.Range("AE5").Value = 36
$ws.Range("AO5").Value = 28
$ws.Range("G5").Value = 2.24
$ws.Range("I5").Value = 3.5
$ws.Range("L5").Value = 1.35
$ws.Range("N5").Value = 4.7
$ws.Range("V5").Value = 1.4
$ws.Range("AI6").Value = 90
$ws.Range("AL6").Value = 42
$ws.Range("F6").Value = 1.82
$ws.Range("S6").Value = 4
$ws.Range("T6").Value = 2.02
$ws.Range("U6").Value = 1.92
$ws.Range("F7").Value = 3.3
$ws.Range("G7").Value = 3.35
$ws.Range("P7").Value = 1.98
$ws.Range("W7").Value = 1.42
$ws.Range("AM8").Value = 110
$ws.Range("Q8").Value = 2.14
$ws.Range("Y8").Value = 15.5
$ws.Range("AE9").Value = 21
$ws.Range("AL9").Value = 46
$ws.Range("S9").Value = 3.1
$ws.Range("U9").Value = 2.32
$ws.Range("AB10").Value = 8.6
$ws.Range("AF10").Value = 14.5
$ws.Range("AH10").Value = 21
$ws.Range("AN10").Value = 32
$ws.Range("AO10").Value = 55
$ws.Range("G10").Value = 2.6
$ws.Range("J10").Value = 3.1
$ws.Range("L10").Value = 1.52
$ws.Range("M10").Value = 1.12
$ws.Range("N10").Value = 3
$ws.Range("O10").Value = 1.47
$ws.Range("P10").Value = 1.67
$ws.Range("Q10").Value = 2.44
$ws.Range("S10").Value = 4.7
$ws.Range("U10").Value = 1.93
$ws.Range("W10").Value = 1.62
$ws.Range("X10").Value = 9.4
$ws.Range("AE11").Value = 55
$ws.Range("F11").Value = 2.42
$ws.Range("G11").Value = 2.44
$ws.Range("H11").Value = 3.65
$ws.Range("L11").Value = 1.52
$ws.Range("N11").Value = 2.98
$ws.Range("P11").Value = 1.66
$ws.Range("Q11").Value = 2.42
$ws.Range("T11").Value = 2
$ws.Range("U11").Value = 1.94
$ws.Range("W11").Value = 1.69
$ws.Range("AJ12").Value = 340
$ws.Range("H12").Value = 1.39
$ws.Range("J12").Value = 5.4
$ws.Range("K12").Value = 5.5
$ws.Range("R12").Value = 1.54
$ws.Range("T12").Value = 2
$ws.Range("V12").Value = 3.45
$ws.Range("Y12").Value = 9.4
$ws.Range("AD13").Value = 20
$ws.Range("AO13").Value = 70
$ws.Range("H13").Value = 5.5
$ws.Range("I13").Value = 5.6
$ws.Range("J13").Value = 4.1
$ws.Range("K13").Value = 4.2
$ws.Range("L13").Value = 1.37
$ws.Range("V13").Value = 1.21
$ws.Range("AM14").Value = 95
$ws.Range("F14").Value = 5.4
$ws.Range("G14").Value = 5.5
$ws.Range("H14").Value = 1.72
$ws.Range("I14").Value = 1.73
$ws.Range("L14").Value = 1.34
$ws.Range("M14").Value = 1.06
$ws.Range("O14").Value = 1.26
$ws.Range("U14").Value = 2.2
$ws.Range("V14").Value = 2.36
$ws.Range("X14").Value = 17.5
